$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<and the force which on hundred four he donged aftentuon, which of oncheod attentuon, when hundred fourse tentuoned it us tentuon.>"
$ws.Range("C2").Value = 60.42742202813713

$ws.Range("B3").Value = "<frired in and acceptent atter day mother work who his regulate was regulate was regulate was regulate was regulate here work cared in accepared in accepe.>"
$ws.Range("C3").Value = 62.51138893838248

$ws.Range("B4").Value = "<in read period appearates of two hundred appearious.>"
$ws.Range("C4").Value = 62.34586634894488

$ws.Range("B5").Value = "<bund to lead an the corner and waited in the transporner, and waited in the transpor.>"
$ws.Range("C5").Value = 58.8772917440528

$ws.Range("B6").Value = "<and to was see of the assued to windowing any discosed the country untry untry untry untry untry untry.>"
$ws.Range("C6").Value = 62.33804337662157

$ws.Range("B7").Value = "<he dine the roded to get curpare to get curpared to get curpare to get curpartant, that he wanted average to get curtant.>"
$ws.Range("C7").Value = 62.08114972593403

$ws.Range("B8").Value = "<while as the act the heales the act the healed the heldor.>"
$ws.Range("C8").Value = 58.80464659166002

$ws.Range("B9").Value = "<toward a fishild a fishild a fishild a fishilds.>"
$ws.Range("C9").Value = 62.38433746850028

$ws.Range("B10").Value = "<in the morinas of marchings he actrongest docket and strong conviction, as to the cond.>"
$ws.Range("C10").Value = 60.34201096545785

$ws.Range("B11").Value = "<tame and every the government to the bitter disardish aimed to the government to the government, every the government,>"
$ws.Range("C11").Value = 62.19708362936724

$ws.Range("B12").Value = "<but it way with themselves and accomplished their perpomplished them away with their perpomplished them away with their perper perpose.>"
$ws.Range("C12").Value = 60.13753306241427

$ws.Range("B13").Value = "<in the came suflic prompted murder rers of convicted murderer rear rear rears almate.>"
$ws.Range("C13").Value = 63.25219145450502

$ws.Range("B14").Value = "<for three is poute secution excution excution, condempt to death an forts poor to death an fort to death an fort to death nemed execution.>"
$ws.Range("C14").Value = 64.61443434788315

$ws.Range("B15").Value = "<agreem due indure to the inagemed to the in agreed states, thing greturn to the in agreem duestabls.>"
$ws.Range("C15").Value = 60.69027380282159

$ws.Range("B16").Value = "<he was stear peoppeant member of the low for appeacted member of thurt people good works.>"
$ws.Range("C16").Value = 60.9611644106381

$ws.Range("B17").Value = "<main dorceiman dorceiman dorsiminan dorsemencey fore.>"
$ws.Range("C17").Value = 60.74827815995827

$ws.Range("B18").Value = "<howere explosied shart serril but shart shart shart shart shart word>"
$ws.Range("C18").Value = 64.60625068688483

$ws.Range("B19").Value = "<he he had a pars a pars at lugggement to a pars at lugggement to a pars at lugggement to a pares at luggement to a palant at lught rockers at lutte.>"
$ws.Range("C19").Value = 62.70361906468717

$ws.Range("B20").Value = "<and instencel of the councels of the secory of the councely if the councely if the councely atterate,>"
$ws.Range("C20").Value = 58.74415505166075

$ws.Range("B21").Value = "<sommons of these deperations of these deperations of the great sill some of the walls moties of the great silless of the great some of the walls.>"
$ws.Range("C21").Value = 63.80927497677936

$ws.Range("B22").Value = "<oswald was all was all was all was all was all was all was all was all was all was all was all.>"
$ws.Range("C22").Value = 64.35565127488536

$ws.Range("B23").Value = "<the service performents preferrst have have have have tuggents prefersence.>"
$ws.Range("C23").Value = 65.91221113357545

$ws.Range("B24").Value = "<howen plained down plained down plained down plained down plained down and sistements of a good agails the require.>"
$ws.Range("C24").Value = 61.26490681500704

$ws.Range("B25").Value = "<the game was near was near was near was near was near was near was near was near was near was near was near was met was metch.>"
$ws.Range("C25").Value = 62.17219159663571

$ws.Range("B26").Value = "<when he had stick get the words with the words with the wordst with the wordstords.>"
$ws.Range("C26").Value = 59.21496596058812

$ws.Range("B27").Value = "<as would nineteen sixty three, nineteen sixty three, nineteen sixty>"
$ws.Range("C27").Value = 60.66672350087008

$ws.Range("B28").Value = "<five points and sixth point six points and sixth point.>"
$ws.Range("C28").Value = 59.879146270214

$ws.Range("B29").Value = "<roval vold revolval vold rover,>"
$ws.Range("C29").Value = 61.81713469875542

$ws.Range("B30").Value = "<lest him against him against him against him against him against him against.>"
$ws.Range("C30").Value = 65.74494181563551

$ws.Range("B31").Value = "<this common the specially appon the specially apponied strops and the mare specially apponied strops and the mare specially apon the specially apponied,>"
$ws.Range("C31").Value = 62.13856230680842

$ws.Range("B32").Value = "<in regarding the insposions to immotory by the commissions to insposed in structionst in struction>"
$ws.Range("C32").Value = 63.88379678191954

$ws.Range("B33").Value = "<two husband shoulds member husbands member husband shoulds member.>"
$ws.Range("C33").Value = 60.83548724638507

$ws.Range("B34").Value = "<who dintravily of the presidents car at the speak who dintravily of the president#s car at the preside.>"
$ws.Range("C34").Value = 61.48208317696258

$ws.Range("B35").Value = "<his appearable for for for for few muld was sentempt sent three propulation one coment one sentement one sent three propulation one signed one signed.>"
$ws.Range("C35").Value = 60.42353352346789

$ws.Range("B36").Value = "<advibe oswald was unnoten such such such such such such such such such sercols.>"
$ws.Range("C36").Value = 58.62539858038343

$ws.Range("B37").Value = "<and secret which would new arrea realed the secret and secret and secret and secret and secret and secret and sevent and state, end quote.>"
$ws.Range("C37").Value = 63.7115049026671

$ws.Range("B38").Value = "<quite pounds arest sums underty pounds a resumber the forms, afthough in ritis.>"
$ws.Range("C38").Value = 61.50959611065412

$ws.Range("B39").Value = "<to have thoughtens by four he to have thought to have at thought to have thousand pounds by froughs.>"
$ws.Range("C39").Value = 62.08074693595204

$ws.Range("B40").Value = "<and told him a stroory told him a story told himmise.>"
$ws.Range("C40").Value = 60.2370059729684

$ws.Range("B41").Value = "<which hard to him and quote, quote, quote,>"
$ws.Range("C41").Value = 63.68324218171881

$ws.Range("B42").Value = "<they tought addus an a romor more more more more more more more more more more.>"
$ws.Range("C42").Value = 57.17872683235013
